$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-03-22", "overview", "K02000001", "United Kingdom", 4301925, 5342, 17, 126172),
    @("2021-03-23", "overview", "K02000001", "United Kingdom", 4307304, 5379, 112, 126284),
    @("2021-03-24", "overview", "K02000001", "United Kingdom", 4312908, 5605, 98, 126382)
)

$startRow = 223
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}
